# Apply the corrected GHI input values (lat/lon fix + recomputed sunrise/sunset
# and clear/cloudy sky figures) to the "Daily" and "Hourly" sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Daily")
$ws2 = $wb.Worksheets.Item("Hourly")

# --- New sunrise / sunset timestamps (shared across every data row) ---
$newSunrise = "2024-02-24T07:11:21"
$newSunset  = "2024-02-24T17:58:05"

# --- Daily sheet : single data row (row 2) ---
$ws1.Range("A2").Value = 47.2229
$ws1.Range("B2").Value = 24.7244
$ws1.Range("E2").Value = $newSunrise
$ws1.Range("F2").Value = $newSunset
$ws1.Range("G2").Value = 3481.73
$ws1.Range("H2").Value = 6806.91
$ws1.Range("I2").Value = 805.79
$ws1.Range("J2").Value = 957.35
$ws1.Range("L2").Value = 957.35

# --- Hourly sheet : rows 2-25, one row per hour of the day ---
for ($r = 2; $r -le 25; $r++) {
    $ws2.Range("A$r").Value = 47.2229
    $ws2.Range("B$r").Value = 24.7244
    $ws2.Range("E$r").Value = $newSunrise
    $ws2.Range("F$r").Value = $newSunset
}

# Daylight hours (rows 9-19) also get refreshed clear/cloudy sky figures.
$ws2.Range("H9").Value = 24.99
$ws2.Range("I9").Value = 151.52
$ws2.Range("J9").Value = 19.69
$ws2.Range("K9").Value = 12.03
$ws2.Range("M9").Value = 12.03

$ws2.Range("H10").Value = 156.62
$ws2.Range("I10").Value = 519.91
$ws2.Range("J10").Value = 58.26
$ws2.Range("K10").Value = 62.55
$ws2.Range("M10").Value = 62.55

$ws2.Range("H11").Value = 305.41
$ws2.Range("I11").Value = 683.35
$ws2.Range("J11").Value = 78.73
$ws2.Range("K11").Value = 91.11
$ws2.Range("M11").Value = 91.11

$ws2.Range("H12").Value = 429.28
$ws2.Range("I12").Value = 766.39
$ws2.Range("J12").Value = 91.04
$ws2.Range("K12").Value = 120.44
$ws2.Range("M12").Value = 120.44

$ws2.Range("H13").Value = 511.44
$ws2.Range("I13").Value = 808.53
$ws2.Range("J13").Value = 97.93
$ws2.Range("K13").Value = 138.96
$ws2.Range("M13").Value = 138.96

$ws2.Range("H14").Value = 542.69
$ws2.Range("I14").Value = 822.69
$ws2.Range("J14").Value = 100.34
$ws2.Range("K14").Value = 144.51
$ws2.Range("M14").Value = 144.51

$ws2.Range("H15").Value = 519.78
$ws2.Range("I15").Value = 812.44
$ws2.Range("J15").Value = 98.57
$ws2.Range("K15").Value = 136.69
$ws2.Range("M15").Value = 136.69

$ws2.Range("H16").Value = 445.09
$ws2.Range("I16").Value = 775.2
$ws2.Range("J16").Value = 92.42
$ws2.Range("K16").Value = 114.44
$ws2.Range("M16").Value = 114.44

$ws2.Range("H17").Value = 326.88
$ws2.Range("I17").Value = 700.09
$ws2.Range("J17").Value = 81.06
$ws2.Range("K17").Value = 81.72
$ws2.Range("M17").Value = 81.72

$ws2.Range("H18").Value = 180.47
$ws2.Range("I18").Value = 554.46
$ws2.Range("J18").Value = 62.22
$ws2.Range("K18").Value = 45.12
$ws2.Range("M18").Value = 45.12

$ws2.Range("H19").Value = 39.08
$ws2.Range("I19").Value = 212.33
$ws2.Range("J19").Value = 25.53
$ws2.Range("K19").Value = 9.77
$ws2.Range("M19").Value = 9.77
